$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label
$ws.Range("B1").Value = "AvgPrice"

# Update individual prices with average prices
$ws.Range("B2").Value = 4240.1
$ws.Range("B3").Value = 2315.4
$ws.Range("B4").Value = 679.4
$ws.Range("B5").Value = 259
$ws.Range("B6").Value = 1082.4
$ws.Range("B7").Value = 502.4

# Update total (sum of average prices)
$ws.Range("B8").Value = 9078.699999999999
